$wb = $excel.ActiveWorkbook

# Productdata sheet: AverageDemand for product 1 (row 2, column G) goes from 40 to 70
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 70

# Work around a round-trip quirk where empty shared-string cells (column H,
# the unused "Name" column) get serialized with a spurious default value on
# save. Re-asserting them as empty keeps them blank as in the source file.
$wsProductdata.Range("H2:H11").Value = ""

# ForecastedAverageDemand sheet: last three periods (rows 9-11, column B) go from 0 to 100
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# ForcastedStandardDeviation sheet: last three periods (rows 9-11, column B) go from 0 to the new values
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
